# Auto-generated edit script applying the cryptos.xlsx data refresh described in the commit diff.
# Each row updates the Price (D) and/or Volume(1h) (E) columns; several rows also had their
# Coin/Link (B/C) swapped because the underlying ranking reordered adjacent coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.904.94"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.633.64"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'214.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "'19.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.860.17"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.663.80"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "'4.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'62.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0754"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "25.904.96"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "'192.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "'9.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'6.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "'143.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "'0.126"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("D28").Value = "'6.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'15.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'0.0497"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("D36").Value = "'0.901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "1.137.29"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'0.548"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'99.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'5.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("D45").Value = "1.769.64"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'56.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.416"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
